$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 5599
$ws1.Range("F8").Value = 6468
$ws1.Range("F13").Value = 47

# Sheet "全部类型" (fourth sheet) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 5599
$ws4.Range("F9").Value = 6468
$ws4.Range("F14").Value = 47
